$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ChambreHip")
$c = $ws.Cells.Item(2,2)
$c.NumberFormat = "@"
$c.Font.ColorIndex = -4105   # xlColorIndexAutomatic
$c.Value = "00010"
Write-Host "B2: " $c.Value()
